$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the species-observation data in rows 4-7 down by one row,
# with row 7's data wrapping around to row 4 (columns A, B, E, F, G, H, L, Q, R).
# Capture the "before" values for each row first, since writes must not
# clobber values we still need to read for later rows.

$rows = 4, 5, 6, 7
$data = @{}

foreach ($r in $rows) {
    $data[$r] = @{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        E = $ws.Cells.Item($r, 5).Value()
        F = $ws.Cells.Item($r, 6).Value()
        G = $ws.Cells.Item($r, 7).Value()
        H = $ws.Cells.Item($r, 8).Value()
        L = $ws.Cells.Item($r, 12).Value()
        Q = $ws.Cells.Item($r, 17).Value()
        R = $ws.Cells.Item($r, 18).Value()
    }
}

# Map: destination row -> source row (row 4 gets row 7's data, row 5 gets
# row 4's data, row 6 gets row 5's data, row 7 gets row 6's data).
$srcFor = @{ 4 = 7; 5 = 4; 6 = 5; 7 = 6 }

foreach ($dst in $rows) {
    $src = $srcFor[$dst]
    $vals = $data[$src]

    $ws.Cells.Item($dst, 1).Value = $vals.A
    $ws.Cells.Item($dst, 2).Value = $vals.B
    $ws.Cells.Item($dst, 5).Value = $vals.E
    $ws.Cells.Item($dst, 6).Value = $vals.F
    $ws.Cells.Item($dst, 7).Value = $vals.G
    $ws.Cells.Item($dst, 8).Value = $vals.H
    $ws.Cells.Item($dst, 17).Value = $vals.Q
    $ws.Cells.Item($dst, 18).Value = $vals.R
}

# Column L ("Kön") holds a blank placeholder cell on some rows and is simply
# absent on others; that presence/absence rotates with the rest of the row
# data. Row 6 currently has the blank placeholder and row 4 needs to gain
# one (sourced, in rotation order, from row 7 -> row 4). Copy the existing
# blank cell (preserving its exact blank representation) before clearing
# the now-stale source, rather than writing/clearing a value directly
# (which deletes the cell outright instead of leaving an empty one behind).
$ws.Cells.Item(6, 12).Copy($ws.Cells.Item(4, 12))
$ws.Cells.Item(6, 12).ClearContents()
